$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix shared string values: correct the misspelling "Golliras" -> "Gorillas"
# and make sure rows keep their correct animal labels.
$ws.Range("A3").Value = "Gorillas"
$ws.Range("A4").Value = "Elephants"

# Update selection to reflect the new active cell A3
$ws.Range("A3").Select()
